$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.31250333333333
$ws.Range("H2").Value = 36.93751
$ws.Range("I2").Value = 0.6498350963072504
$ws.Range("J2").Value = 0.6498350963072506
$ws.Range("M2").Value = 45.924193
$ws.Range("N2").Value = 137.772579
$ws.Range("O2").Value = 0.307792367338991
$ws.Range("P2").Value = 0.307792367338991
$ws.Range("Q2").Value = 565.4417793931434
$ws.Range("R2").Value = 5088.976014538291
$ws.Range("S2").Value = 0.2000142826723698
$ws.Range("T2").Value = 0.2000142826723699
$ws.Range("G3").Value = 12.31250333333333
$ws.Range("H3").Value = 36.93751
$ws.Range("I3").Value = 0.6498350963072504
$ws.Range("J3").Value = 0.6498350963072506
$ws.Range("O3").Value = 0.2696759485354523
$ws.Range("P3").Value = 0.2696759485354523
$ws.Range("Q3").Value = 495.4185495817622
$ws.Range("R3").Value = 4458.766946235861
$ws.Range("S3").Value = 0.1752448959882847
$ws.Range("T3").Value = 0.1752448959882848
$ws.Range("G4").Value = 12.31250333333333
$ws.Range("H4").Value = 36.93751
$ws.Range("I4").Value = 0.6498350963072504
$ws.Range("J4").Value = 0.6498350963072506
$ws.Range("M4").Value = 23.60320766666667
$ws.Range("N4").Value = 70.809623
$ws.Range("O4").Value = 0.1581930283351338
$ws.Range("P4").Value = 0.1581930283351339
$ws.Range("Q4").Value = 290.6145730731922
$ws.Range("R4").Value = 2615.53115765873
$ws.Range("S4").Value = 0.1027993818032973
$ws.Range("T4").Value = 0.1027993818032973
$ws.Range("G5").Value = 12.31250333333333
$ws.Range("H5").Value = 36.93751
$ws.Range("I5").Value = 0.6498350963072504
$ws.Range("J5").Value = 0.6498350963072506
$ws.Range("M5").Value = 19.226538
$ws.Range("N5").Value = 57.679614
$ws.Range("O5").Value = 0.1288597852280838
$ws.Range("P5").Value = 0.1288597852280838
$ws.Range("Q5").Value = 236.72681321346
$ws.Range("R5").Value = 2130.54131892114
$ws.Range("S5").Value = 0.08373761094382345
$ws.Range("T5").Value = 0.08373761094382347
$ws.Range("G6").Value = 12.31250333333333
$ws.Range("H6").Value = 36.93751
$ws.Range("I6").Value = 0.6498350963072504
$ws.Range("J6").Value = 0.6498350963072506
$ws.Range("M6").Value = 20.21413933333333
$ws.Range("N6").Value = 60.64241799999999
$ws.Range("O6").Value = 0.1354788705623391
$ws.Range("P6").Value = 0.1354788705623391
$ws.Range("Q6").Value = 248.8866579221311
$ws.Range("R6").Value = 2239.97992129918
$ws.Range("S6").Value = 0.08803892489947514
$ws.Range("T6").Value = 0.08803892489947515
$ws.Range("I7").Value = 0.3333514949915254
$ws.Range("J7").Value = 0.3333514949915254
$ws.Range("M7").Value = 45.924193
$ws.Range("N7").Value = 137.772579
$ws.Range("O7").Value = 0.307792367338991
$ws.Range("P7").Value = 0.307792367338991
$ws.Range("Q7").Value = 290.0595298137787
$ws.Range("R7").Value = 2610.535768324008
$ws.Range("S7").Value = 0.1026030457994334
$ws.Range("T7").Value = 0.1026030457994334
$ws.Range("I8").Value = 0.3333514949915254
$ws.Range("J8").Value = 0.3333514949915254
$ws.Range("O8").Value = 0.2696759485354523
$ws.Range("P8").Value = 0.2696759485354523
$ws.Range("S8").Value = 0.08989688060755069
$ws.Range("T8").Value = 0.0898968806075507
$ws.Range("I9").Value = 0.3333514949915254
$ws.Range("J9").Value = 0.3333514949915254
$ws.Range("M9").Value = 23.60320766666667
$ws.Range("N9").Value = 70.809623
$ws.Range("O9").Value = 0.1581930283351338
$ws.Range("P9").Value = 0.1581930283351339
$ws.Range("Q9").Value = 149.0790555185218
$ws.Range("R9").Value = 1341.711499666696
$ws.Range("S9").Value = 0.0527338824927536
$ws.Range("T9").Value = 0.05273388249275361
$ws.Range("I10").Value = 0.3333514949915254
$ws.Range("J10").Value = 0.3333514949915254
$ws.Range("M10").Value = 19.226538
$ws.Range("N10").Value = 57.679614
$ws.Range("O10").Value = 0.1288597852280838
$ws.Range("P10").Value = 0.1288597852280838
$ws.Range("Q10").Value = 121.435788152592
$ws.Range("R10").Value = 1092.922093373328
$ws.Range("S10").Value = 0.04295560205006862
$ws.Range("T10").Value = 0.04295560205006863
$ws.Range("I11").Value = 0.3333514949915254
$ws.Range("J11").Value = 0.3333514949915254
$ws.Range("M11").Value = 20.21413933333333
$ws.Range("N11").Value = 60.64241799999999
$ws.Range("O11").Value = 0.1354788705623391
$ws.Range("P11").Value = 0.1354788705623391
$ws.Range("Q11").Value = 127.6735282123929
$ws.Range("R11").Value = 1149.061753911536
$ws.Range("S11").Value = 0.0451620840417191
$ws.Range("T11").Value = 0.04516208404171911
$ws.Range("G12").Value = 0.3185656666666667
$ws.Range("H12").Value = 0.955697
$ws.Range("I12").Value = 0.01681340870122405
$ws.Range("J12").Value = 0.01681340870122405
$ws.Range("M12").Value = 45.924193
$ws.Range("N12").Value = 137.772579
$ws.Range("O12").Value = 0.307792367338991
$ws.Range("P12").Value = 0.307792367338991
$ws.Range("Q12").Value = 14.62987115917367
$ws.Range("R12").Value = 131.668840432563
$ws.Range("S12").Value = 0.005175038867187741
$ws.Range("T12").Value = 0.005175038867187741
$ws.Range("G13").Value = 0.3185656666666667
$ws.Range("H13").Value = 0.955697
$ws.Range("I13").Value = 0.01681340870122405
$ws.Range("J13").Value = 0.01681340870122405
$ws.Range("O13").Value = 0.2696759485354523
$ws.Range("P13").Value = 0.2696759485354523
$ws.Range("Q13").Value = 12.81813586188245
$ws.Range("R13").Value = 115.363222756942
$ws.Range("S13").Value = 0.004534171939616823
$ws.Range("T13").Value = 0.004534171939616823
$ws.Range("G14").Value = 0.3185656666666667
$ws.Range("H14").Value = 0.955697
$ws.Range("I14").Value = 0.01681340870122405
$ws.Range("J14").Value = 0.01681340870122405
$ws.Range("M14").Value = 23.60320766666667
$ws.Range("N14").Value = 70.809623
$ws.Range("O14").Value = 0.1581930283351338
$ws.Range("P14").Value = 0.1581930283351339
$ws.Range("Q14").Value = 7.519171585803445
$ws.Range("R14").Value = 67.67254427223101
$ws.Range("S14").Value = 0.002659764039082922
$ws.Range("T14").Value = 0.002659764039082922
$ws.Range("G15").Value = 0.3185656666666667
$ws.Range("H15").Value = 0.955697
$ws.Range("I15").Value = 0.01681340870122405
$ws.Range("J15").Value = 0.01681340870122405
$ws.Range("M15").Value = 19.226538
$ws.Range("N15").Value = 57.679614
$ws.Range("O15").Value = 0.1288597852280838
$ws.Range("P15").Value = 0.1288597852280838
$ws.Range("Q15").Value = 6.124914895662001
$ws.Range("R15").Value = 55.124234060958
$ws.Range("S15").Value = 0.002166572234191727
$ws.Range("T15").Value = 0.002166572234191727
$ws.Range("G16").Value = 0.3185656666666667
$ws.Range("H16").Value = 0.955697
$ws.Range("I16").Value = 0.01681340870122405
$ws.Range("J16").Value = 0.01681340870122405
$ws.Range("M16").Value = 20.21413933333333
$ws.Range("N16").Value = 60.64241799999999
$ws.Range("O16").Value = 0.1354788705623391
$ws.Range("P16").Value = 0.1354788705623391
$ws.Range("Q16").Value = 6.439530772816222
$ws.Range("R16").Value = 57.955776955346
$ws.Range("S16").Value = 0.002277861621144839
$ws.Range("T16").Value = 0.002277861621144839
